# SCD0018-006 - Update TC_ID from "DGS-298" to "SCD0018-006"
# (commit: "Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the new TC_ID family (SCD0283 -> SCD0018)
$ws.Name = "SCD0018"

# Update the TC_ID column (B) for every test-step row from the old "DGS-298"
# value to the new "SCD0018-006" value.
$ws.Range("B2").Value = "SCD0018-006"
$ws.Range("B3").Value = "SCD0018-006"
$ws.Range("B4").Value = "SCD0018-006"

# The longer TC_ID text no longer fits the old best-fit column width / row
# height, so widen column B and grow row 2 to fit the new content.
$ws.Columns.Item(2).ColumnWidth = 11.67
$ws.Rows.Item(2).RowHeight = 89.25

# Move the active selection from D4 to B5.
$ws.Range("B5").Select()
